$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill column U (U1:U20) with values 0..19
for ($r = 1; $r -le 20; $r++) {
    $ws.Cells.Item($r, 21).Value = $r - 1
}

# Fill new row 21 (A21:T21) with values 0..19
for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item(21, $c).Value = $c - 1
}

# Update the selected cell to match the final state of the diff
$ws.Range("Y19").Select()
